# feat: add 2022-Q4 data
#
# 1. Insert a new worksheet "2022-Q4" right after "总计", pushing the other
#    quarter sheets (2022-Q3, 2022-Q1, 2020-Q4) one slot to the right.
# 2. Populate "2022-Q4" with the same fund-holding table shape used by the
#    other quarter sheets.
# 3. Insert a new row into "总计" (right after the header) summarizing the
#    2022-Q4 quarter, shifting the previously existing summary rows down.

$wb = $excel.ActiveWorkbook

$total = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# Step 1: create + place the new "2022-Q4" worksheet right after "总计"
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Add($null, $total)
$q4.Name = "2022-Q4"

# Match the outline / page-setup bits other sheets carry (<sheetPr>).
$q4.Outline.SummaryRow = 1
$q4.Outline.SummaryColumn = 1

# Match the page margins used by the sibling quarter sheets.
$q4.PageSetup.LeftMargin = 54
$q4.PageSetup.RightMargin = 54
$q4.PageSetup.TopMargin = 72
$q4.PageSetup.BottomMargin = 72
$q4.PageSetup.HeaderMargin = 36
$q4.PageSetup.FooterMargin = 36

# Header row (bold, centered, boxed - same look as the other quarter sheets).
$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $q4.Cells.Item(1, 2 + $c).Value = $headers[$c]
}
# Re-use the exact header style already used by the sibling quarter sheets
# (copy its format rather than rebuilding font/border/alignment by hand, so
# the cells land on the very same style index).
$total.Range("B1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)

# Fund rows.
$rows = @(
    @("0", "501030", "汇添富中证环境治理指数（LOF）A", "2.99", "92.14", "1.92", "0.0574", 9),
    @("1", "164908", "交银施罗德中证环境治理指数（LOF）", "1.55", "93.92", "1.97", "0.0305", 9),
    @("2", "501031", "汇添富中证环境治理指数（LOF）C", "1.38", "92.14", "1.92", "0.0265", 9),
    @("3", "013413", "交银施罗德中证环境治理指数（LOF）C", "0.11", "93.92", "1.97", "0.0022", 9)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = 2 + $i
    $data = $rows[$i]

    $q4.Cells.Item($r, 1).Value = [int]$data[0]

    # Columns B..G hold text that merely looks numeric (e.g. "501030",
    # "2.99") in the source data - force text storage, then drop the
    # number-format override so no stray style index is left behind.
    $textRange = $q4.Range($q4.Cells.Item($r, 2), $q4.Cells.Item($r, 7))
    $textRange.NumberFormat = "@"
    $q4.Cells.Item($r, 2).Value = $data[1]
    $q4.Cells.Item($r, 3).Value = $data[2]
    $q4.Cells.Item($r, 4).Value = $data[3]
    $q4.Cells.Item($r, 5).Value = $data[4]
    $q4.Cells.Item($r, 6).Value = $data[5]
    $q4.Cells.Item($r, 7).Value = $data[6]
    $textRange.ClearFormats()

    $q4.Cells.Item($r, 8).Value = $data[7]
}

# Copy the "A" column's bold/centered/boxed style (same one used by every
# other quarter sheet's row-number column) onto A2:A5.
$total.Range("A2").Copy()
$q4.Range("A2:A5").PasteSpecial(-4122)
for ($i = 0; $i -lt $rows.Length; $i++) {
    $q4.Cells.Item(2 + $i, 1).Value = [int]$rows[$i][0]
}

# ---------------------------------------------------------------------
# Step 2: add the 2022-Q4 summary row to "总计" (row 2, pushing the rest
# down, matching the order already used for the other quarters).
# ---------------------------------------------------------------------
$total.Rows.Item(2).Insert()
$total.Range("B2:D2").ClearFormats()

# Re-use the (still correctly styled) style that now sits on A3 - it is an
# exact copy of what every other "A" summary cell in this column uses.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 4
$total.Range("D2").Value = 0.12

# ---------------------------------------------------------------------
# Step 3: restore the originally-active sheet (2020-Q4, the last tab) so
# the new sheet doesn't steal the "active" / tabSelected state.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Activate()
